$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 436.66666
$ws.Range("I12").Value = 405
$ws.Range("K12").Value = 405
$ws.Range("M12").Value = -235
$ws.Range("H17").Value = 2912.0667
$ws.Range("J17").Value = 2912.0667
$ws.Range("L17").Value = 8736.2001
$ws.Range("N17").Value = -9072.2001
$ws.Range("H33").Value = 458.5
$ws.Range("I33").Value = 458.5
$ws.Range("K33").Value = 458.5
$ws.Range("M33").Value = -229.5
$ws.Range("H39").Value = 139.33333
$ws.Range("I39").Value = 139.33333
$ws.Range("K39").Value = 417.99999
$ws.Range("M39").Value = -121.99999
$ws.Range("H58").Value = 3606.3333
$ws.Range("I58").Value = 401
$ws.Range("J58").Value = 10017
$ws.Range("K58").Value = 1203
$ws.Range("L58").Value = 30051
$ws.Range("M58").Value = -1053
$ws.Range("N58").Value = -30351
$ws.Range("H92").Value = 728.4
$ws.Range("I92").Value = 584.5454999999999
$ws.Range("J92").Value = 1124
$ws.Range("K92").Value = 584.5454999999999
$ws.Range("L92").Value = 1124
$ws.Range("M92").Value = 663.4545000000001
$ws.Range("N92").Value = -3620
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H137").Value = 1998.8154
$ws.Range("I137").Value = 1460.1875
$ws.Range("J137").Value = 3519.647
$ws.Range("K137").Value = 4380.5625
$ws.Range("L137").Value = 10558.941
$ws.Range("M137").Value = -1830.5625
$ws.Range("N137").Value = -15658.941
$ws.Range("H138").Value = 5071.579
$ws.Range("J138").Value = 5178.9375
$ws.Range("L138").Value = 15536.8125
$ws.Range("N138").Value = -25816.8125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2912
$ws.Range("J2").Value = 2830.3333
$ws.Range("L2").Value = 2830.3333
$ws.Range("N2").Value = -3056.3333
$ws.Range("H74").Value = 2013.7273
$ws.Range("I74").Value = 1316.8
$ws.Range("K74").Value = 1316.8
$ws.Range("M74").Value = -442.8
$ws.Range("H77").Value = 2013.7273
$ws.Range("I77").Value = 1316.8
$ws.Range("K77").Value = 6584
$ws.Range("M77").Value = -2216
$ws.Range("H97").Value = 249
$ws.Range("I97").Value = 190.5
$ws.Range("K97").Value = 190.5
$ws.Range("M97").Value = 305.5
$ws.Range("H116").Value = 2912
$ws.Range("J116").Value = 2830.3333
$ws.Range("L116").Value = 2830.3333
$ws.Range("N116").Value = -7418.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2912
$ws.Range("J3").Value = 2830.3333
$ws.Range("L3").Value = 2830.3333
$ws.Range("N3").Value = -3058.3333
$ws.Range("H105").Value = 7620.222
$ws.Range("I105").Value = 7654.7144
$ws.Range("J105").Value = 7499.5
$ws.Range("K105").Value = 7654.7144
$ws.Range("L105").Value = 7499.5
$ws.Range("M105").Value = -5907.7144
$ws.Range("N105").Value = -10993.5
$ws.Range("H107").Value = 1105
$ws.Range("I107").Value = 947.2353000000001
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 947.2353000000001
$ws.Range("L107").Value = 1999
$ws.Range("M107").Value = 972.7646999999999
$ws.Range("N107").Value = -5839
$ws.Range("H134").Value = 3995.182
$ws.Range("I134").Value = 3993.75
$ws.Range("K134").Value = 11981.25
$ws.Range("M134").Value = -9446.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4590.3335
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 4385.5
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 4385.5
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -5633.5
$ws.Range("H65").Value = 4590.3335
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 4385.5
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 21927.5
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -28167.5
$ws.Range("H107").Value = 1499.8823
$ws.Range("I107").Value = 990.1
$ws.Range("K107").Value = 990.1
$ws.Range("M107").Value = 929.9
$ws.Range("H141").Value = 62967.94
$ws.Range("J141").Value = 62967.94
$ws.Range("L141").Value = 62967.94
$ws.Range("N141").Value = -73327.94

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 168499.5
$ws.Range("J122").Value = 251999.25
$ws.Range("L122").Value = 2267993.25
$ws.Range("N122").Value = -2272893.25
$ws.Range("H127").Value = 750
$ws.Range("I127").Value = 500
$ws.Range("K127").Value = 1500
$ws.Range("M127").Value = 3460

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1762.7778
$ws.Range("I132").Value = 1702.625
$ws.Range("J132").Value = 2244
$ws.Range("K132").Value = 5107.875
$ws.Range("L132").Value = 6732
$ws.Range("M132").Value = -2577.875
$ws.Range("N132").Value = -11792

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2341.7144
$ws.Range("I68").Value = 1932.3334
$ws.Range("J68").Value = 2648.75
$ws.Range("K68").Value = 1932.3334
$ws.Range("L68").Value = 2648.75
$ws.Range("M68").Value = -1183.3334
$ws.Range("N68").Value = -4146.75
$ws.Range("H71").Value = 2341.7144
$ws.Range("I71").Value = 1932.3334
$ws.Range("J71").Value = 2648.75
$ws.Range("K71").Value = 9661.666999999999
$ws.Range("L71").Value = 13243.75
$ws.Range("M71").Value = -5917.666999999999
$ws.Range("N71").Value = -20731.75
$ws.Range("H82").Value = 2247.8572
$ws.Range("I82").Value = 2214.6
$ws.Range("J82").Value = 2331
$ws.Range("K82").Value = 2214.6
$ws.Range("L82").Value = 2331
$ws.Range("M82").Value = -1853.6
$ws.Range("N82").Value = -3053
$ws.Range("H85").Value = 2247.8572
$ws.Range("I85").Value = 2214.6
$ws.Range("J85").Value = 2331
$ws.Range("K85").Value = 2214.6
$ws.Range("L85").Value = 2331
$ws.Range("M85").Value = -966.5999999999999
$ws.Range("N85").Value = -4827
$ws.Range("H93").Value = 2300.2778
$ws.Range("I93").Value = 2369.3333
$ws.Range("K93").Value = 2369.3333
$ws.Range("M93").Value = -1121.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 986.3333
$ws.Range("I96").Value = 986.3333
$ws.Range("K96").Value = 986.3333
$ws.Range("M96").Value = 386.6667
$ws.Range("H132").Value = 3360.875
$ws.Range("I132").Value = 2617.0908
$ws.Range("K132").Value = 7851.2724
$ws.Range("M132").Value = -5321.2724
